# Commit: "updated xpath for billingpage"
#
# 1. Rename the "CreatingOrders" sheet to "BillingCreateOrders".
# 2. Update its selection from A1:C1 to a single-cell selection on E1.
# 3. Make the "Billing" sheet the active tab and set its selection to G1
#    (this moves tabSelected + activeTab off of "AddingAsset" and onto "Billing").

$wb = $excel.ActiveWorkbook

# Rename CreatingOrders -> BillingCreateOrders
$wsOrders = $wb.Worksheets.Item("CreatingOrders")
$wsOrders.Name = "BillingCreateOrders"

# Update selection on the renamed sheet: A1:C1 -> E1
$wsOrders.Range("E1").Select()

# Activate Billing sheet and update its selection: A1:C1 -> G1
$wsBilling = $wb.Worksheets.Item("Billing")
$wsBilling.Activate()
$wsBilling.Range("G1").Select()
